# Update the "Metadata" sheet of the CDA Logical Model workbook
# (StructureDefinition-PQR) to reflect the new publication:
#   - Version bumped:        2.0.0-sd-202406-matchbox-patch -> 2.0.1-sd-202510-matchbox-patch
#   - Date bumped:            2024-06-19T17:47:42+02:00 -> 2025-10-29T22:15:57+01:00
#   - A new "Jurisdiction" property row is inserted right after "Contact"
#     (with an empty value), pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update Version value (row 3, column B) ---
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# --- Update Date value (row 8, column B) ---
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Insert a new row after "Contact" (row 10) for "Jurisdiction" ---
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row above (Contact, row 10) onto the freshly
# inserted row 11 so it matches the rest of the table (border/alignment).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
